$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Worksheet")

# Update "Contenu du stage" (programming language) statistics: rows 16-23
# Column D = language name, E = number of students, G = percentage text.
# G values look like percentages, so force them to stay as literal text
# (matching the source file, where they are plain strings, not number-
# formatted cells) by temporarily applying a text format and then
# clearing the format again so no style index is left on the cell.

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# C# : 34 -> 13 students, 97.14 % -> 37.14 %
$ws.Range("E16").Value = 13
Set-TextValue $ws.Range("G16") "37.14 %"

# COBOL : 0 -> 15 students, 0 % -> 42.86 %
$ws.Range("E17").Value = 15
Set-TextValue $ws.Range("G17") "42.86 %"

# C++ : stays 0 students, 0 %
$ws.Range("E18").Value = 0
Set-TextValue $ws.Range("G18") "0 %"

# ASSEMBLEUR : 1 -> 2 students, 2.86 % -> 5.71 %
$ws.Range("E19").Value = 2
Set-TextValue $ws.Range("G19") "5.71 %"

# ANDROID : 0 -> 5 students, 0 % -> 14.29 %
$ws.Range("E20").Value = 5
Set-TextValue $ws.Range("G20") "14.29 %"

# JEE : stays 0 students, 0 %
$ws.Range("E21").Value = 0
Set-TextValue $ws.Range("G21") "0 %"

# DELPHI : stays 0 students, 0 %
$ws.Range("E22").Value = 0
Set-TextValue $ws.Range("G22") "0 %"

# PHP5 : stays 0 students, 0 %
$ws.Range("E23").Value = 0
Set-TextValue $ws.Range("G23") "0 %"
